$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# --- Insert 7 new blank rows right after row 48 (rows 49-55), shifting the
# --- existing blank/bottom-border rows down by 7. Row 48 itself (already an
# --- empty, pre-formatted entry row) is left untouched.
$ws.Rows("49:55").Insert()

# The Insert() call resets the brand-new rows to each column's default
# style instead of the "empty log entry" style used throughout the table
# (medium box border on B..E/G, date/time number formats). Restore that
# formatting by copying it from row 56 (the first row below the inserted
# block, which still carries the original formatting untouched).
$ws.Range("B56:E56").Copy()
$ws.Range("B49:E55").PasteSpecial(-4122)
$ws.Range("G56").Copy()
$ws.Range("G49:G55").PasteSpecial(-4122)
$ws.Rows("49:55").RowHeight = 24.95

# --- Fill in the four new activity-log entries (row 48, which already
# --- existed, plus the first three freshly inserted rows 49-51).
$ws.Range("B48").Value = 6977
$ws.Range("C48").Value = 43926
$ws.Range("D48").Value = 0.60416666666666663
$ws.Range("E48").Value = 0.61805555555555558
$ws.Range("G48").Value = "Provided tech support to give team members access to Github repository"

$ws.Range("B49").Value = 6977
$ws.Range("C49").Value = 43926
$ws.Range("D49").Value = 0.61805555555555558
$ws.Range("E49").Value = 0.63541666666666663
$ws.Range("G49").Value = "Reviewed team mates' work on Arith.vhd, Adder.vhd and related files"

$ws.Range("B50").Value = 6977
$ws.Range("C50").Value = 43926
$ws.Range("D50").Value = 0.63541666666666663
$ws.Range("E50").Value = 0.64236111111111105
$ws.Range("G50").Value = "Fixed formatting of Arith.vhd and Adder.vhd."

$ws.Range("A51").Value = "s"
$ws.Range("B51").Value = 6977
$ws.Range("C51").Value = 43926
$ws.Range("D51").Value = 0.64236111111111105
$ws.Range("G51").Value = "Removed unnecessary VHDL code from LogicUnit.vhd."

# --- Match the author's final selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C51").Select()

Write-Host "Applied activity-log updates."
